$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated S-val stat data (regenerated to filter save games)
# Each row: row#, B=TB, C=d2S, D=K, E=IP, G=sum (F=Win unchanged)
$data = @(
    @(2, 1.459612070389937, 10.29869402782916, 0.8054896365839992, 8.660232485948974, 21.22402822075207),
    @(3, 1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(4, 1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(5, 0.01514828764759746, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1.35982162114495),
    @(6, 0.04763786555579896, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 12.65070551845588),
    @(7, 1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 29.84159230404497),
    @(8, 1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286),
    @(9, 1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759),
    @(10, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(11, 1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759),
    @(12, 3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671),
    @(13, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(14, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(15, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(16, 1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(17, 0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.642425054193055),
    @(18, 1.459612070389937, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 2.80428540388729),
    @(19, 1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759),
    @(20, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(21, 1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 29.84159230404497),
    @(22, 1.459612070389937, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 3.074671312995807),
    @(23, 1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 29.84159230404497),
    @(24, 0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 2.290389397800092),
    @(25, 1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 29.84159230404497),
    @(26, 0.127881588408715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.837881874639075),
    @(27, 0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044),
    @(28, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(29, 3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(30, 1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(31, 3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(32, 1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 12.59312877619104)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B - TB
    $ws.Cells.Item($r, 3).Value = $row[2]   # C - d2S
    $ws.Cells.Item($r, 4).Value = $row[3]   # D - K
    $ws.Cells.Item($r, 5).Value = $row[4]   # E - IP
    $ws.Cells.Item($r, 7).Value = $row[5]   # G - sum
}
